$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data for row 2 has the values in columns A/B, C/D, and E/F swapped
# (pairwise transposition), per the diff:
#   A2 <-> B2
#   C2 <-> D2
#   E2 <-> F2
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2

$ws.Range("A2").Value = $b2
$ws.Range("B2").Value = $a2
$ws.Range("C2").Value = $d2
$ws.Range("D2").Value = $c2
$ws.Range("E2").Value = $f2
$ws.Range("F2").Value = $e2

# Reflect the new active-cell selection recorded in the saved view state.
$ws.Range("G20").Select() | Out-Null
